# Applies the "Creado diseño de creacion de artwork (no funcional, faltan ajustar SELECTS)" edit
# to PlaTreball.xlsx: updates several "Durada real" (I column) estimate values on Sprint2,
# and fills in two previously-blank task rows (32 and 33) with new task data that references
# two newly introduced descriptions.

$wb = $excel.ActiveWorkbook

# Sprint2 is the second worksheet (tabSelected sheet) in the workbook.
$ws = $wb.Worksheets.Item(2)

# --- Update "Durada real" (column I) values for existing tasks ---
$ws.Cells.Item(16, 9).Value = 2      # I16: 3 -> 2
$ws.Cells.Item(18, 9).Value = 2      # I18: 3 -> 2
$ws.Cells.Item(19, 9).Value = 1.5    # I19: 3 -> 1.5
$ws.Cells.Item(20, 9).Value = 3.5    # I20: 3 -> 3.5
$ws.Cells.Item(21, 9).Value = 2.5    # I21: 3 -> 2.5
$ws.Cells.Item(22, 9).Value = 2      # I22: 3 -> 2
$ws.Cells.Item(24, 9).Value = 2      # I24: (blank) -> 2

# --- Fill in previously empty row 32 with a new task ---
$ws.Cells.Item(32, 3).Value = "Marc Martin"
$ws.Cells.Item(32, 4).Value = "marcmartin60"
$ws.Cells.Item(32, 5).Value = "Creació de progress bar + disseny"
$ws.Cells.Item(32, 6).Value = "Leonard Craciun"
$ws.Cells.Item(32, 7).Value = 2
$ws.Cells.Item(32, 8).Value = 1
$ws.Cells.Item(32, 9).Value = 1

# --- Fill in previously empty row 33 with a new task ---
$ws.Cells.Item(33, 3).Value = "Marc Martin"
$ws.Cells.Item(33, 4).Value = "marcmartin60"
$ws.Cells.Item(33, 5).Value = "Añadida opció de filtrar per estat amb camps de la bdd"
$ws.Cells.Item(33, 6).Value = "Leonard Craciun"
$ws.Cells.Item(33, 7).Value = 0.5
$ws.Cells.Item(33, 8).Value = 1
$ws.Cells.Item(33, 9).Value = 0.4

# --- Update the view state of Sprint2 so the window is scrolled to the new rows ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D40").Select() | Out-Null
